# Applies the cryptos list refresh described in the commit:
# "Updated cryptos list on Wed Nov 22 08:45:01 UTC 2023 with GitHub Actions"
#
# Each entry updates one cell. "AsText" forces the NumberFormat to Text
# before the write (then resets the style) so that numeric-looking price
# strings such as "58.30" are not auto-coerced into the Double 58.3 by
# Excel's normal smart-typing, matching the source data which stores
# prices/volumes as plain text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "36.495.10"; AsText = $true },
    @{ Cell = "E2"; Value = "  -2.45%  "; AsText = $false },
    @{ Cell = "D3"; Value = "2.000.90"; AsText = $true },
    @{ Cell = "E3"; Value = "  -0.39%  "; AsText = $false },
    @{ Cell = "E4"; Value = "  -0.03%  "; AsText = $false },
    @{ Cell = "D5"; Value = "235.04"; AsText = $true },
    @{ Cell = "E5"; Value = "  -9.27%  "; AsText = $false },
    @{ Cell = "E6"; Value = "  -3.17%  "; AsText = $false },
    @{ Cell = "E7"; Value = "  +0.01%  "; AsText = $false },
    @{ Cell = "D8"; Value = "54.73"; AsText = $true },
    @{ Cell = "E8"; Value = "  -2.16%  "; AsText = $false },
    @{ Cell = "E9"; Value = "  -3.53%  "; AsText = $false },
    @{ Cell = "D10"; Value = "58.30"; AsText = $true },
    @{ Cell = "E10"; Value = "  +3.34%  "; AsText = $false },
    @{ Cell = "E11"; Value = "  -2.84%  "; AsText = $false },
    @{ Cell = "D12"; Value = "0.0980"; AsText = $true },
    @{ Cell = "E12"; Value = "  -3.68%  "; AsText = $false },
    @{ Cell = "D13"; Value = "2.294.73"; AsText = $true },
    @{ Cell = "E13"; Value = "  -0.39%  "; AsText = $false },
    @{ Cell = "E14"; Value = "  -0.30%  "; AsText = $false },
    @{ Cell = "D15"; Value = "20.30"; AsText = $true },
    @{ Cell = "E15"; Value = "  -2.89%  "; AsText = $false },
    @{ Cell = "E16"; Value = "  -5.48%  "; AsText = $false },
    @{ Cell = "D17"; Value = "5.06"; AsText = $true },
    @{ Cell = "E17"; Value = "  -3.14%  "; AsText = $false },
    @{ Cell = "D18"; Value = "2.006.93"; AsText = $true },
    @{ Cell = "E18"; Value = "  -0.04%  "; AsText = $false },
    @{ Cell = "D19"; Value = "36.444.18"; AsText = $true },
    @{ Cell = "E19"; Value = "  -2.20%  "; AsText = $false },
    @{ Cell = "D20"; Value = "67.78"; AsText = $true },
    @{ Cell = "E20"; Value = "  -2.71%  "; AsText = $false },
    @{ Cell = "D21"; Value = "0.0₃0804"; AsText = $true },
    @{ Cell = "E21"; Value = "  -3.77%  "; AsText = $false },
    @{ Cell = "D22"; Value = "5.28"; AsText = $true },
    @{ Cell = "E22"; Value = "  +2.84%  "; AsText = $false },
    @{ Cell = "D23"; Value = "221.45"; AsText = $true },
    @{ Cell = "E23"; Value = "  -2.97%  "; AsText = $false },
    @{ Cell = "E24"; Value = "  +0.09%  "; AsText = $false },
    @{ Cell = "D25"; Value = "2.36"; AsText = $true },
    @{ Cell = "E25"; Value = "  +1.10%  "; AsText = $false },
    @{ Cell = "E26"; Value = "  -8.24%  "; AsText = $false },
    @{ Cell = "D27"; Value = "162.86"; AsText = $true },
    @{ Cell = "E27"; Value = "  -1.10%  "; AsText = $false },
    @{ Cell = "E28"; Value = "  -2.25%  "; AsText = $false },
    @{ Cell = "D29"; Value = "1.37"; AsText = $true },
    @{ Cell = "E29"; Value = "  +3.16%  "; AsText = $false },
    @{ Cell = "E30"; Value = "  -3.55%  "; AsText = $false },
    @{ Cell = "D31"; Value = "18.70"; AsText = $true },
    @{ Cell = "E31"; Value = "  -4.55%  "; AsText = $false },
    @{ Cell = "E32"; Value = "  -2.66%  "; AsText = $false },
    @{ Cell = "D33"; Value = "4.37"; AsText = $true },
    @{ Cell = "E33"; Value = "  -4.98%  "; AsText = $false },
    @{ Cell = "D34"; Value = "0.0604"; AsText = $true },
    @{ Cell = "E34"; Value = "  -6.07%  "; AsText = $false },
    @{ Cell = "B35"; Value = "LidoDAOToken"; AsText = $false },
    @{ Cell = "C35"; Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; AsText = $false },
    @{ Cell = "D35"; Value = "2.38"; AsText = $true },
    @{ Cell = "E35"; Value = "  +0.97%  "; AsText = $false },
    @{ Cell = "B36"; Value = "InternetComputer(DFINITY)"; AsText = $false },
    @{ Cell = "C36"; Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; AsText = $false },
    @{ Cell = "D36"; Value = "4.26"; AsText = $true },
    @{ Cell = "E36"; Value = "  -5.55%  "; AsText = $false },
    @{ Cell = "E37"; Value = "  -0.05%  "; AsText = $false },
    @{ Cell = "B38"; Value = "RenderToken"; AsText = $false },
    @{ Cell = "C38"; Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; AsText = $false },
    @{ Cell = "D38"; Value = "3.31"; AsText = $true },
    @{ Cell = "E38"; Value = "  +0.00%  "; AsText = $false },
    @{ Cell = "B39"; Value = "WEMIXToken"; AsText = $false },
    @{ Cell = "C39"; Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; AsText = $false },
    @{ Cell = "D39"; Value = "1.76"; AsText = $true },
    @{ Cell = "E39"; Value = "  -3.18%  "; AsText = $false },
    @{ Cell = "D40"; Value = "5.62"; AsText = $true },
    @{ Cell = "E40"; Value = "  +5.25%  "; AsText = $false },
    @{ Cell = "E41"; Value = "  -1.65%  "; AsText = $false },
    @{ Cell = "D42"; Value = "1.452.36"; AsText = $true },
    @{ Cell = "E42"; Value = "  +2.64%  "; AsText = $false },
    @{ Cell = "D43"; Value = "0.0927"; AsText = $true },
    @{ Cell = "E43"; Value = "  +0.26%  "; AsText = $false },
    @{ Cell = "E44"; Value = "  -4.55%  "; AsText = $false },
    @{ Cell = "D45"; Value = "89.54"; AsText = $true },
    @{ Cell = "E45"; Value = "  +0.40%  "; AsText = $false },
    @{ Cell = "E46"; Value = "  -7.99%  "; AsText = $false },
    @{ Cell = "D47"; Value = "15.16"; AsText = $true },
    @{ Cell = "E47"; Value = "  -3.48%  "; AsText = $false },
    @{ Cell = "D48"; Value = "0.994"; AsText = $true },
    @{ Cell = "E48"; Value = "  -2.63%  "; AsText = $false },
    @{ Cell = "D49"; Value = "2.88"; AsText = $true },
    @{ Cell = "E49"; Value = "  -1.12%  "; AsText = $false },
    @{ Cell = "B50"; Value = "SynthetixNetwork"; AsText = $false },
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"; AsText = $false },
    @{ Cell = "D50"; Value = "3.75"; AsText = $true },
    @{ Cell = "E50"; Value = "  +7.48%  "; AsText = $false },
    @{ Cell = "B51"; Value = "FraxShare"; AsText = $false },
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; AsText = $false },
    @{ Cell = "D51"; Value = "6.85"; AsText = $true },
    @{ Cell = "E51"; Value = "  -2.26%  "; AsText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.AsText) {
        # Pin the cell to Text format so the numeric-looking string survives
        # verbatim (trailing zeros, thousand-separator dots, etc.), then drop
        # back to the workbook's default "Normal" style so no stray
        # explicit number format is left behind on the cell.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
